# Update view-count figures (column F) on the "展览", "演出" and "全部类型"
# sheets to match the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2686
$ws1.Range("F5").Value = 1498
$ws1.Range("F6").Value = 1136
$ws1.Range("F11").Value = 117
$ws1.Range("F13").Value = 9135
$ws1.Range("F14").Value = 394
$ws1.Range("F15").Value = 2498
$ws1.Range("F16").Value = 4
$ws1.Range("F17").Value = 257
$ws1.Range("F18").Value = 181
$ws1.Range("F20").Value = 628
$ws1.Range("F23").Value = 999
$ws1.Range("F24").Value = 2085
$ws1.Range("F25").Value = 2175
$ws1.Range("F27").Value = 1882
$ws1.Range("F28").Value = 1927
$ws1.Range("F30").Value = 1463
$ws1.Range("F31").Value = 272
$ws1.Range("F33").Value = 207
$ws1.Range("F37").Value = 289
$ws1.Range("F38").Value = 487
$ws1.Range("F40").Value = 29
$ws1.Range("F41").Value = 419
$ws1.Range("F42").Value = 14
$ws1.Range("F43").Value = 1379
$ws1.Range("F45").Value = 7
$ws1.Range("F46").Value = 8
$ws1.Range("F47").Value = 615
$ws1.Range("F49").Value = 296

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 21

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2686
$ws4.Range("F4").Value = 1498
$ws4.Range("F6").Value = 1136
$ws4.Range("F9").Value = 117
$ws4.Range("F10").Value = 9135
$ws4.Range("F11").Value = 394
$ws4.Range("F12").Value = 2498
$ws4.Range("F14").Value = 4
$ws4.Range("F16").Value = 257
$ws4.Range("F17").Value = 181
$ws4.Range("F19").Value = 628
$ws4.Range("F21").Value = 999
$ws4.Range("F22").Value = 2175
$ws4.Range("F23").Value = 1882
$ws4.Range("F25").Value = 1463
$ws4.Range("F26").Value = 272
$ws4.Range("F28").Value = 207
$ws4.Range("F32").Value = 289
$ws4.Range("F33").Value = 487
$ws4.Range("F34").Value = 21
$ws4.Range("F38").Value = 29
$ws4.Range("F39").Value = 419
$ws4.Range("F41").Value = 14
$ws4.Range("F42").Value = 1379
$ws4.Range("F45").Value = 7
$ws4.Range("F46").Value = 8
$ws4.Range("F47").Value = 615
$ws4.Range("F48").Value = 296
